$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for new column F (adds a new shared string) ---
$ws.Range("F1").Value2 = "Swapped Ifelse order H_ext = 2"

# --- Column C updates: rows whose value flips from 0 to 1 ---
$cRows = @(3, 4, 15, 16, 20, 25, 26, 28, 29, 30, 33, 37, 41, 42, 45, 49, 50, 51, 52, 54, 55, 57, 58, 60, 63, 64, 65, 66, 69, 70, 71, 72, 73, 76, 77, 78, 80, 81, 83, 84, 86, 88, 89, 92, 98, 100, 101)
foreach ($r in $cRows) {
    $ws.Cells.Item($r, 3).Value2 = 1
}

# --- New boolean column F, rows 2-101 ---
$fVals = @{
    2 = $true
    3 = $false
    4 = $true
    5 = $true
    6 = $false
    7 = $false
    8 = $false
    9 = $false
    10 = $false
    11 = $true
    12 = $false
    13 = $false
    14 = $false
    15 = $false
    16 = $false
    17 = $false
    18 = $false
    19 = $false
    20 = $false
    21 = $true
    22 = $false
    23 = $false
    24 = $false
    25 = $false
    26 = $true
    27 = $false
    28 = $true
    29 = $false
    30 = $false
    31 = $false
    32 = $false
    33 = $false
    34 = $false
    35 = $true
    36 = $false
    37 = $true
    38 = $false
    39 = $false
    40 = $false
    41 = $false
    42 = $true
    43 = $false
    44 = $false
    45 = $false
    46 = $false
    47 = $false
    48 = $true
    49 = $true
    50 = $true
    51 = $true
    52 = $true
    53 = $true
    54 = $true
    55 = $true
    56 = $false
    57 = $true
    58 = $false
    59 = $true
    60 = $false
    61 = $false
    62 = $true
    63 = $false
    64 = $true
    65 = $true
    66 = $true
    67 = $true
    68 = $false
    69 = $false
    70 = $true
    71 = $true
    72 = $false
    73 = $false
    74 = $true
    75 = $false
    76 = $false
    77 = $true
    78 = $false
    79 = $false
    80 = $false
    81 = $true
    82 = $false
    83 = $false
    84 = $true
    85 = $false
    86 = $false
    87 = $true
    88 = $false
    89 = $false
    90 = $false
    91 = $true
    92 = $false
    93 = $true
    94 = $false
    95 = $true
    96 = $true
    97 = $false
    98 = $false
    99 = $true
    100 = $true
    101 = $true
}
foreach ($r in $fVals.Keys) {
    $ws.Cells.Item($r, 6).Value2 = $fVals[$r]
}

# --- Column F width ---
$ws.Columns.Item(6).ColumnWidth = 24.3

# --- Page setup (portrait orientation, 1200 dpi print quality) ---
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PrintQuality = 1200

# --- View / selection changes ---
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B51").Select() | Out-Null
